$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add I0 and IF headers in columns I and J, matching the
# existing header style (copy format from H1, the last header cell).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for new columns I (I0) and J (IF)
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 3

$ws.Range("I4").Value = 6
$ws.Range("J4").Value = 6

$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 3

$ws.Range("I6").Value = 8
$ws.Range("J6").Value = 8

$ws.Range("I7").Value = 6
$ws.Range("J7").Value = 8

$ws.Range("I8").Value = 11
$ws.Range("J8").Value = 11

$ws.Range("I9").Value = 6
$ws.Range("J9").Value = 7

$ws.Range("I10").Value = 7
$ws.Range("J10").Value = 8
